$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted right before the current row 160,
# pushing every following row (160-282) down by one (161-283).
$ws.Rows("160:160").Insert()

# Populate the newly inserted row 160 with the new price record.
# Most fields mirror the (now shifted-down) neighbouring record, except
# for the reporting date (D) and the volume (J), which are new values.
$ws.Cells.Item(160, 1).Value = 8
$ws.Cells.Item(160, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(160, 3).Value = "Coquimbo"
$ws.Cells.Item(160, 4).Value = 44729
$ws.Cells.Item(160, 5).Value = 4
$ws.Cells.Item(160, 6).Value = 100112003
$ws.Cells.Item(160, 7).Value = "Ajo"
$ws.Cells.Item(160, 8).Value = "Chino"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 560
$ws.Cells.Item(160, 11).Value = 19000
$ws.Cells.Item(160, 12).Value = 20000
$ws.Cells.Item(160, 13).Value = 19500
$ws.Cells.Item(160, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(160, 15).Value = "China"
$ws.Cells.Item(160, 16).Value = 1950
$ws.Cells.Item(160, 17).Value = 10
$ws.Cells.Item(160, 18).Value = "Hortaliza"
